# Sort the data table by the "type" column (B), keeping a stable order
# for rows that share the same type (matches "add sort for type" commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRow = 6
$firstDataRow = 7
$lastDataRow = 67
$lastCol = 8  # column H

$dataRange = $ws.Range($ws.Cells.Item($firstDataRow, 2), $ws.Cells.Item($lastDataRow, $lastCol))
$keyRange = $ws.Range($ws.Cells.Item($firstDataRow, 2), $ws.Cells.Item($lastDataRow, 2))

$dataRange.Sort($keyRange, 1)
